# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @{
    2  = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 6.48142807727062, 12.0302756157461)
    3  = @(1.505614041169197, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 6.741336633845642)
    4  = @(1.505614041169197, 1.65323645889881, 0.7127328510149897, 6.48142807727062, 10.35301142835362)
    5  = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 6.048734245549538)
    6  = @(3.182878228561681, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 5.488907176552729)
    7  = @(3.182878228561681, 1.65323645889881, 0.1529057820181812, 6.48142807727062, 11.47044854674929)
    8  = @(1.505614041169197, 1.65323645889881, 16.98373111632243, 0.4998867070740569, 20.64246832346449)
    9  = @(3.182878228561681, 9.226618575922256, 16.98373111632243, 6.48142807727062, 35.87465599807698)
    10 = @(3.182878228561681, 9.226618575922256, 3.082599426703578, 246.9852506941017, 262.4773469252892)
    11 = @(0.3464964993005633, 9.226618575922256, 3.082599426703578, 246.9852506941017, 259.640965196028)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals[0]
    $ws.Range("C$row").Value = $vals[1]
    $ws.Range("D$row").Value = $vals[2]
    $ws.Range("E$row").Value = $vals[3]
    $ws.Range("G$row").Value = $vals[4]
}
